$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.BottomPadding = 9999999
